$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row correct-answer mark value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row correct marks value (B12): 42 -> 70
$ws.Range("B12").Value = 70

# Update the correct/total marks text (E12): "36/84" -> "70/140"
$ws.Range("E12").Value = "70/140"
